# Add "Sheet2" right after "Sheet1" and populate it with the API route
# reference table described in the commit diff.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a brand-new sheet immediately after Sheet1 (Before=$null, After=$ws1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ----- header row -----
$ws2.Range("D1").Value = "메소드"
$ws2.Range("E1").Value = "url"
$ws2.Range("F1").Value = "파라미터"

# ----- 홈 -----
$ws2.Range("B2").Value = "홈"
$ws2.Range("C2").Value = "게시물전체정보"
$ws2.Range("D2").Value = "get"
$ws2.Range("E2").Value = "/home"

# ----- 나침반 -----
$ws2.Range("B3").Value = "나침반"
$ws2.Range("C3").Value = "게시물전체정보"
$ws2.Range("D3").Value = "get"
$ws2.Range("E3").Value = "/main"

# ----- 게시물 (merged B4:B8) -----
$ws2.Range("B4").Value = "게시물"
$ws2.Range("C4").Value = "특정게시물정보"
$ws2.Range("D4").Value = "get"
$ws2.Range("E4").Value = "/post/{id}"

$ws2.Range("C5").Value = "저장"
$ws2.Range("D5").Value = "post"
$ws2.Range("E5").Value = "/post"
$ws2.Range("F5").Value = "id"

$ws2.Range("C6").Value = "삭제"
$ws2.Range("D6").Value = "post"
$ws2.Range("E6").Value = "/post/delete"
$ws2.Range("F6").Value = "id"

$ws2.Range("C7").Value = "수정"
$ws2.Range("D7").Value = "post"
$ws2.Range("E7").Value = "/post/update"
$ws2.Range("F7").Value = "id"

$ws2.Range("C8").Value = "좋아요"
$ws2.Range("D8").Value = "post"
$ws2.Range("E8").Value = "/post/like/"
$ws2.Range("F8").Value = "id"

# ----- 댓글 (merged B9:B11) -----
$ws2.Range("B9").Value = "댓글"
$ws2.Range("C9").Value = "저장"
$ws2.Range("D9").Value = "post"
$ws2.Range("E9").Value = "/comment"
$ws2.Range("F9").Value = "id"

$ws2.Range("C10").Value = "삭제"
$ws2.Range("D10").Value = "post"
$ws2.Range("E10").Value = "/comment/delete"
$ws2.Range("F10").Value = "id"

$ws2.Range("C11").Value = "좋아요"
$ws2.Range("D11").Value = "post"
$ws2.Range("E11").Value = "/comment/like"
$ws2.Range("F11").Value = "id"

# ----- footer note -----
$ws2.Range("B13").Value = "페이징기법"
$ws2.Range("C13").Value = "보류"

# ----- merges -----
$ws2.Range("B4:B8").Merge()
$ws2.Range("B9:B11").Merge()

# ----- alignment (center / center) on the node labels -----
$centered = $ws2.Range("B2:B3,B4:B11,B12")
$centered.HorizontalAlignment = -4108
$centered.VerticalAlignment = -4108

# ----- column widths (bestFit-style, closest reachable values) -----
$ws2.Columns.Item(2).ColumnWidth = 10.2857142857143
$ws2.Columns.Item(3).ColumnWidth = 14.4285714285714
$ws2.Columns.Item(5).ColumnWidth = 16.1428571428571

# ----- selection / active sheet -----
$ws2.Activate()
$ws2.Range("E13").Select()
